$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
# F1 needs the same direct formatting (bold / border / centered / top)
# that the other header cells (B1:E1) already carry, so copy that
# formatting from E1 before we repurpose E1's text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# F1 becomes the old "average debt per person" header; E1 becomes the
# new "% of total people" header.
$ws.Range("F1").Value = "average debt per person"
$ws.Range("E1").Value = "% of total people"

# --- Data rows -----------------------------------------------------------
# The occupation rows are reordered (merchant now first, then gentleman,
# then weaver), a new "% of total people" column is inserted at E, and the
# old "average debt per person" figures move out to the new column F.

# Row 2 -> merchant (was row 3's data)
$ws.Range("B2").Value = "merchant"
$ws.Range("C2").Value = 90.99623299999999
$ws.Range("D2").Value = 3
$ws.Range("F2").Value = 30.33207766666666
$ws.Range("E2").Value = 60

# Row 3 -> gentleman (was row 2's data)
$ws.Range("B3").Value = "gentleman"
$ws.Range("C3").Value = 78.33
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = 78.33
$ws.Range("E3").Value = 20

# Row 4 -> weaver (unchanged occupation/debt total/# of people; only the
# new % column is added and the average debt figure shifts to F)
$ws.Range("F4").Value = 46.02
$ws.Range("E4").Value = 20
